$wb = $excel.ActiveWorkbook

# --- Add sheet "t3" right after "t2", populate with actual/device readings ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$t3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$t3.Name = "t3"
$t3.Range("A1").Value = "actual"
$t3.Range("B1").Value = "device"
$t3.Range("A2").Value = 36.1
$t3.Range("B2").Value = 33.88
$t3.Range("A3").Value = 36.2
$t3.Range("B3").Value = 34.5
$t3.Range("A4").Value = 36.2
$t3.Range("B4").Value = 34.31
$t3.Range("A5").Value = 36.3
$t3.Range("B5").Value = 34.31
$t3.Range("A6").Value = 36.2
$t3.Range("B6").Value = 34.31
$t3.Range("A7").Value = 36.2
$t3.Range("B7").Value = 34.13
$t3.Range("A8").Value = 36.2
$t3.Range("B8").Value = 34.25
$t3.Range("A9").Value = 36.4
$t3.Range("B9").Value = 34.56
$t3.Range("A10").Value = 36.3
$t3.Range("B10").Value = 34.63
$t3.Range("A11").Value = 36.1
$t3.Range("B11").Value = 34.56
$t3.Range("A12").Value = 36
$t3.Range("B12").Value = 34.56
$t3.Range("A13").Value = 36.2
$t3.Range("B13").Value = 34.44
$t3.Range("A14").Value = 36.2
$t3.Range("B14").Value = 34.63
$t3.Range("A15").Value = 36.1
$t3.Range("B15").Value = 34.56
$t3.Range("A16").Value = 36.1
$t3.Range("B16").Value = 34.5
[void]$t3.Range("B17").Select()

# --- Add sheet "t4" right after "t3", populate with actual/device readings ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$t4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$t4.Name = "t4"
$t4.Range("A1").Value = "actual"
$t4.Range("B1").Value = "device"
$t4.Range("A2").Value = 35.7
$t4.Range("B2").Value = 34.31
$t4.Range("A3").Value = 36.3
$t4.Range("B3").Value = 34.25
$t4.Range("A4").Value = 36.3
$t4.Range("B4").Value = 34.19
$t4.Range("A5").Value = 36.1
$t4.Range("B5").Value = 34.56
$t4.Range("A6").Value = 36
$t4.Range("B6").Value = 34.06
$t4.Range("A7").Value = 36.1
$t4.Range("A8").Value = 36.3
$t4.Range("A9").Value = 36.5
$t4.Range("A10").Value = 36.5
$t4.Range("A11").Value = 36.5
$t4.Range("A12").Value = 36.5
$t4.Range("A13").Value = 36.3
$t4.Range("A14").Value = 36.6
$t4.Range("A15").Value = 36.6
$t4.Range("A16").Value = 36.7
$t4.Range("A17").Value = 36.5
$t4.Range("A18").Value = 36.5
$t4.Range("A19").Value = 36.3
$t4.Range("A20").Value = 36.5
$t4.Range("A21").Value = 36.7
$t4.Range("A22").Value = 36.5
$t4.Range("A23").Value = 36.5
$t4.Range("A24").Value = 36.5
$t4.Range("A25").Value = 36.5
$t4.Range("A26").Value = 36.6

# t4 is the final active sheet/tab, scrolled & selected near the bottom of its data
[void]$t4.Select()
[void]$t4.Range("A27").Select()
